# Auto-generated edit script for Products data.xlsx
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Mobile & Computer Accessories
$ws2 = $wb.Worksheets.Item(2)   # Gadgets

### 1) Sheet2 (Gadgets): insert new 'Features head point' column at D, shifting MRP/Sell Price/Product Link right
$ws2.Columns.Item(4).Insert()

### Copy the highlighted data-row style (currently s=6 on Sheet1 D6:D8) onto the new Sheet2 D2:D11 cells
$ws1.Range("D6").Copy()
$ws2.Range("D2:D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

### Header for new column D1
$ws2.Range("D1").Value = "Features head point"

### Feature head point text per Gadget row (D2:D11)
$feat116 = @'
Here are key feature head points for your rechargeable torch:  
- **USB Rechargeable** – Convenient and fast charging with included cable  
- **Ergonomic Grip** – Designed for comfort and easy handling  
- **Bright Illumination** – Powerful beam for clear visibility in the dark  
- **Versatile Use** – Ideal for outdoor activities, emergencies, and daily needs  
- **Durable Build** – Long-lasting and reliable performance  
Let me know if you'd like further refinements!  
'@
$ws2.Range("D2").Value = $feat116

$feat117 = @'
Here are concise and impactful head points for your calculation assistance request:  
- **Math Support** – Get help with basic and complex equations  
- **Arithmetic & Percentages** – Solve everyday numerical problems effortlessly  
- **Problem-Solving** – Assistance with step-by-step calculations  
- **Quick & Accurate** – Reliable solutions tailored to your needs  
- **Versatile Help** – From simple math to advanced computations  
Let me know if you'd like any refinements!  
'@
$ws2.Range("D3").Value = $feat117

$feat118 = @'
Here are sharp and compelling feature head points for your rechargeable torch with a holder:  
- **Powerful Illumination** – Bright beam for exploring, camping & emergencies  
- **Hands-Free Convenience** – Secure holder for easy placement & usability  
- **USB Rechargeable** – Hassle-free charging for reliable performance  
- **Durable Build** – Designed for long-lasting use in any situation  
- **Versatile Companion** – Ideal for outdoor adventures & everyday needs  
Let me know if you’d like any refinements!
'@
$ws2.Range("D4").Value = $feat118

$feat119 = @'
Here are sharp and engaging feature head points for your advanced trimmer:  
- **Precision Grooming** – Achieve effortless styling with accurate trimming  
- **Built-in Digital Meter** – Real-time updates on battery life & performance  
- **Versatile Usage** – Ideal for facial hair, beards & hairstyling needs  
- **Compact & Sleek Design** – User-friendly and travel-ready convenience  
- **Tailored Experience** – Adjustable settings for personalized grooming  
Let me know if you'd like any refinements!
'@
$ws2.Range("D5").Value = $feat119

$feat120 = @'
Here are sharp and compelling feature head points for your rechargeable metal keychain:  
- **Built-in LED Light** – Bright illumination for night-time convenience  
- **Durable Metal Design** – Sleek, lightweight & perfect for everyday carry  
- **USB Rechargeable** – Eco-friendly & battery-saving functionality  
- **Compact & Stylish** – A practical accessory with a modern look  
- **Reliable Utility** – Ideal for unlocking doors & navigating dark spaces  
Let me know if you’d like any refinements!
'@
$ws2.Range("D6").Value = $feat120

$feat121 = @'
Here are impactful feature head points for your handheld megaphone:  
- **Powerful 150W Output** – Ensures clear and loud communication  
- **Built-in Recorder** – Pre-record messages for consistent playback  
- **Ergonomic Handheld Design** – Comfortable grip for easy portability  
- **Ideal for Public Use** – Perfect for announcements, rallies & events  
- **Reliable Sound Amplification** – Projects your voice far & wide  
Let me know if you'd like any refinements!
'@
$ws2.Range("D7").Value = $feat121

$feat122 = @'
Here are concise and engaging feature head points for your rechargeable LED flashlight:  
- **Powerful Beam** – Bright & focused illumination for dark spaces  
- **USB Rechargeable** – No hassle of replacing batteries, always ready  
- **Compact & Lightweight** – Easy to carry for outdoor & emergency use  
- **Durable Design** – Built for longevity & reliable performance  
- **Versatile Utility** – Ideal for adventures, daily tasks & safety  
Let me know if you’d like any refinements!
'@
$ws2.Range("D8").Value = $feat122

$feat123 = @'
Here are impactful feature head points for your scientific calculator:  
- **Advanced Functions** – Trigonometry, logarithms, exponents & statistics  
- **User-Friendly Interface** – Simplifies complex calculations effortlessly  
- **Precision & Efficiency** – Ensures accurate results for math & science  
- **Durable Design** – Built for long-term professional & academic use  
- **Versatile Utility** – Ideal for students, professionals & researchers  
Let me know if you’d like any refinements!
'@
$ws2.Range("D9").Value = $feat123

$feat124 = @'
Here are concise and engaging feature head points for your 8-in-1 board:  
- **Multi-Functionality** – Eight integrated tools for enhanced efficiency  
- **USB Charging Port** – Keep your devices powered with ease  
- **Compact & Durable** – Designed for home, office & travel convenience  
- **Streamlined Productivity** – Ideal for multitasking & organized workflow  
- **Modern & Practical** – The perfect companion for daily use  
Let me know if you’d like any refinements!
'@
$ws2.Range("D10").Value = $feat124

$feat125 = @'
Here are sharp and compelling feature head points for your 3-In-1 Torch:  
- **Powerful Flashlight** – Focused brightness for clear visibility  
- **Wide-Angle Lamp** – Illuminates larger areas for convenience  
- **Emergency Strobe Light** – Essential signaling in critical situations  
- **Rechargeable & Efficient** – Hassle-free power solution for reliability  
- **Durable & Sleek Design** – Built for longevity with modern aesthetics  
- **Versatile Utility** – Ideal for outdoor adventures, emergencies & daily use  
Let me know if you’d like any refinements!
'@
$ws2.Range("D11").Value = $feat125

### New listing-status columns H (under category page) and I (actual product page)
$ws2.Range("H1").Value = "Listing status (under category page)"
$ws2.Range("I1").Value = "Listing status (actual product page)"

### Match header style used on Sheet1 H1:I1 (s=5)
$ws1.Range("H1").Copy()
$ws2.Range("H1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

### Mark every Gadget row as Listed on both listing-status columns
for ($r = 2; $r -le 11; $r++) {
    $ws2.Range("H$r").Value = "Listed"
    $ws2.Range("I$r").Value = "Listed"
}

### Header row 1 on Sheet2 grows taller once it carries wrapped feature-head-point text
$ws2.Rows.Item(1).RowHeight = 45

### 2) Sheet1 (Mobile & Computer Accessories): normalise D6:D8 styling to the plain data style
$ws1.Range("C6").Copy()
$ws1.Range("D6:D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

### Re-apply the original text (PasteSpecial formats-only should not disturb the value, but make sure)
### Fill I9:I14 on Sheet1 with Listed status (previously blank)
for ($r = 9; $r -le 14; $r++) {
    $ws1.Range("I$r").Value = "Listed"
}

### 3) Selection / active-sheet bookkeeping to mirror the saved workbook view
$ws1.Range("H1:I1").Select()
$ws2.Activate()
$ws2.Range("C17").Select()

